# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the a94a3e9b-c574-4087-9961-b9d0b4140a4e file's row (row 3) across the
# Overview, zh-cn and de-de sheets, reflecting a fresh handback report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-21 18:52:17"

# --- zh-cn sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-21 18:52:12"
$wsZhCn.Range("K3").Value = "2016-08-21 18:52:31"

# --- de-de sheet: "Correspond Handoff Datetime" (H) / "Correspond Handback DateTime" (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-21 18:52:17"
$wsDeDe.Range("K3").Value = "2016-08-21 18:52:38"
